$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 136.4330001899934
$ws.Range("B3").Value = 163.6436659911002
$ws.Range("B4").Value = 172.8607842485089
$ws.Range("B5").Value = 181.9214503578388
$ws.Range("B6").Value = 199.785439055213
$ws.Range("B7").Value = 205.9547641225778
$ws.Range("B8").Value = 220.3971529742912
$ws.Range("B9").Value = 235.0497901439441
$ws.Range("B10").Value = 247.5909114315268
$ws.Range("B11").Value = 255.1109172637569
$ws.Range("B12").Value = 262.7569266497218
$ws.Range("B13").Value = 275.8809256278674
$ws.Range("B14").Value = 287.3274738669512
$ws.Range("B15").Value = 299.0575203880709
$ws.Range("B16").Value = 307.5138449752932
$ws.Range("B17").Value = 329.6138629390286
$ws.Range("B18").Value = 339.4019635782203
$ws.Range("B19").Value = 358.1656308301314
$ws.Range("B20").Value = 370.9437000845166
$ws.Range("B21").Value = 372.997120499747
$ws.Range("B22").Value = 377.7280772993946
$ws.Range("B23").Value = 382.3275450050054
$ws.Range("B24").Value = 398.4318941415669
$ws.Range("B25").Value = 413.6587479405925
$ws.Range("B26").Value = 422.0123903175068
$ws.Range("B27").Value = 438.3690659174299
$ws.Range("B28").Value = 449.2243424193895
$ws.Range("B29").Value = 453.702635103083
$ws.Range("B30").Value = 464.3342704125094
$ws.Range("B31").Value = 515.6333318777102
$ws.Range("B32").Value = 515.2474954834963
$ws.Range("B33").Value = 533.688137223595
$ws.Range("B34").Value = 543.352441996271
$ws.Range("B35").Value = 569.1668006838763
$ws.Range("B36").Value = 575.6420470888818
$ws.Range("B37").Value = 585.4507704061054
$ws.Range("B38").Value = 604.4347730903877
$ws.Range("B39").Value = 612.3359978218533
$ws.Range("B40").Value = 629.0964194563225
$ws.Range("B41").Value = 636.6687699528651
$ws.Range("B42").Value = 647.0619079827162
$ws.Range("B43").Value = 651.611472132797
$ws.Range("B44").Value = 658.920603788008
$ws.Range("B45").Value = 680.7530064505553
$ws.Range("B46").Value = 694.2825835001809
$ws.Range("B47").Value = 688.0981489222677
$ws.Range("B48").Value = 705.4050740639814
$ws.Range("B49").Value = 696.8193675717018
$ws.Range("B50").Value = 698.623515216703
$ws.Range("B51").Value = 707.496936660839
$ws.Range("B52").Value = 724.3214951031605
$ws.Range("B53").Value = 727.8953008724416
$ws.Range("B54").Value = 736.4031472035331
$ws.Range("B55").Value = 748.3914581934595
$ws.Range("B56").Value = 745.0841051967426
$ws.Range("B57").Value = 760.0231407208379
$ws.Range("B58").Value = 766.6683821055805
$ws.Range("B59").Value = 768.391092501607
$ws.Range("B60").Value = 771.2482476034353
$ws.Range("B61").Value = 771.1346497471062
$ws.Range("B62").Value = 776.6082943289654
